$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three source rows below were removed in the refreshed export,
# so drop them here too (delete bottom-up to keep row numbers stable).
$ws.Rows.Item(52).Delete()
$ws.Rows.Item(48).Delete()
$ws.Rows.Item(44).Delete()

# Add the two new lookup columns coming from the polygon match: PD (Q) and N2 (R)
$ws.Range("Q1").Value = "PD"
$ws.Range("R1").Value = "N2"
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("Q2").Value = "AGU-E"
$ws.Range("R2").Value = "Fuera de Poligono OVL"
$ws.Range("Q3").Value = "PUE-B"
$ws.Range("R3").Value = "Fuera de Poligono OVL"
$ws.Range("Q4").Value = "BLO-?"
$ws.Range("R4").Value = "Fuera de Poligono OVL"
$ws.Range("Q5").Value = "VCR-I"
$ws.Range("R5").Value = "Fuera de Poligono OVL"
$ws.Range("Q6").Value = "AGU-O"
$ws.Range("R6").Value = "Fuera de Poligono OVL"
$ws.Range("Q7").Value = "BLO-A"
$ws.Range("R7").Value = "Fuera de Poligono OVL"
$ws.Range("Q8").Value = "COG-H"
$ws.Range("R8").Value = "Fuera de Poligono OVL"
$ws.Range("Q9").Value = "DEV-F"
$ws.Range("R9").Value = "ARATO-25058.PO.1DEV"
$ws.Range("Q10").Value = "PUE-G"
$ws.Range("R10").Value = "Fuera de Poligono OVL"
$ws.Range("Q11").Value = "BLO-J"
$ws.Range("R11").Value = "Fuera de Poligono OVL"
$ws.Range("Q12").Value = "BLO-J"
$ws.Range("R12").Value = "Fuera de Poligono OVL"
$ws.Range("Q13").Value = "AGU-N"
$ws.Range("R13").Value = "Fuera de Poligono OVL"
$ws.Range("Q14").Value = "COG-F"
$ws.Range("R14").Value = "Fuera de Poligono OVL"
$ws.Range("Q15").Value = "DEV-L"
$ws.Range("R15").Value = "ARATO-25058.PO.2DEV"
$ws.Range("Q16").Value = "DEV-F"
$ws.Range("R16").Value = "ARATO-25058.PO.1DEV"
$ws.Range("Q17").Value = "VCR-O"
$ws.Range("R17").Value = "Fuera de Poligono OVL"
$ws.Range("Q18").Value = "PUE-M"
$ws.Range("R18").Value = "ARATO-25058.PO.1PUE"
$ws.Range("Q19").Value = "PUE-H"
$ws.Range("R19").Value = "Fuera de Poligono OVL"
$ws.Range("Q20").Value = "ATH-S"
$ws.Range("R20").Value = "Fuera de Poligono OVL"
$ws.Range("Q21").Value = "PUE-M"
$ws.Range("R21").Value = "ARATO-25058.PO.1PUE"
$ws.Range("Q22").Value = "PUE-J"
$ws.Range("R22").Value = "Fuera de Poligono OVL"
$ws.Range("Q23").Value = "PUE-M"
$ws.Range("R23").Value = "ARATO-25058.PO.1PUE"
$ws.Range("Q24").Value = "COG-A"
$ws.Range("R24").Value = "Fuera de Poligono OVL"
$ws.Range("Q25").Value = "PUE-F"
$ws.Range("R25").Value = "Fuera de Poligono OVL"
$ws.Range("Q26").Value = "PUE-I"
$ws.Range("R26").Value = "Fuera de Poligono OVL"
$ws.Range("Q27").Value = "BLO-?"
$ws.Range("R27").Value = "Fuera de Poligono OVL"
$ws.Range("Q28").Value = "ALM-N"
$ws.Range("R28").Value = "Fuera de Poligono OVL"
$ws.Range("Q29").Value = "ALM-O"
$ws.Range("R29").Value = "Fuera de Poligono OVL"
$ws.Range("Q30").Value = "ATH-J"
$ws.Range("R30").Value = "Fuera de Poligono OVL"
$ws.Range("Q31").Value = "CON-B"
$ws.Range("R31").Value = "Fuera de Poligono OVL"
$ws.Range("Q32").Value = "RET-A"
$ws.Range("R32").Value = "Fuera de Poligono OVL"
$ws.Range("Q33").Value = "BLO-N"
$ws.Range("R33").Value = "Fuera de Poligono OVL"
$ws.Range("Q34").Value = "COG-?"
$ws.Range("R34").Value = "Fuera de Poligono OVL"
$ws.Range("Q35").Value = "BLO-F"
$ws.Range("R35").Value = "Fuera de Poligono OVL"
$ws.Range("Q36").Value = "COG-K"
$ws.Range("R36").Value = "Fuera de Poligono OVL"
$ws.Range("Q37").Value = "AGU-C"
$ws.Range("R37").Value = "Fuera de Poligono OVL"
$ws.Range("Q38").Value = "PAV-V"
$ws.Range("R38").Value = "Fuera de Poligono OVL"
$ws.Range("Q39").Value = "ALM-A"
$ws.Range("R39").Value = "Fuera de Poligono OVL"
$ws.Range("Q40").Value = "CON-G"
$ws.Range("R40").Value = "Fuera de Poligono OVL"
$ws.Range("Q41").Value = "PCH-S"
$ws.Range("R41").Value = "ARATO-25058.PO.2PCH"
$ws.Range("Q42").Value = "CLI-I"
$ws.Range("R42").Value = "Fuera de Poligono OVL"
$ws.Range("Q43").Value = "CLI-N"
$ws.Range("R43").Value = "Fuera de Poligono OVL"
$ws.Range("Q44").Value = "AGU-N"
$ws.Range("R44").Value = "Fuera de Poligono OVL"
$ws.Range("Q45").Value = "BLO-I"
$ws.Range("R45").Value = "Fuera de Poligono OVL"
$ws.Range("Q46").Value = "VCR-B"
$ws.Range("R46").Value = "Fuera de Poligono OVL"
$ws.Range("Q47").Value = "ATH-C"
$ws.Range("R47").Value = "Fuera de Poligono OVL"
$ws.Range("Q48").Value = "NRA-R"
$ws.Range("R48").Value = "ARATO-25058.PO.2NRA"
$ws.Range("Q49").Value = "CEN-M"
$ws.Range("R49").Value = "Fuera de Poligono OVL"
$ws.Range("Q50").Value = "NRA-R"
$ws.Range("R50").Value = "ARATO-25058.PO.2NRA"
$ws.Range("Q51").Value = "VCR-M"
$ws.Range("R51").Value = "Fuera de Poligono OVL"
$ws.Range("Q52").Value = "COG-H"
$ws.Range("R52").Value = "Fuera de Poligono OVL"
$ws.Range("Q53").Value = "PAV-M"
$ws.Range("R53").Value = "Fuera de Poligono OVL"
$ws.Range("Q54").Value = "CON-I"
$ws.Range("R54").Value = "Fuera de Poligono OVL"
$ws.Range("Q55").Value = "BLO-F"
$ws.Range("R55").Value = "Fuera de Poligono OVL"
$ws.Range("Q56").Value = "CON-G"
$ws.Range("R56").Value = "Fuera de Poligono OVL"
$ws.Range("Q57").Value = "COG-I"
$ws.Range("R57").Value = "Fuera de Poligono OVL"
$ws.Range("Q58").Value = "CLI-O"
$ws.Range("R58").Value = "Fuera de Poligono OVL"
$ws.Range("Q59").Value = "PUE-J"
$ws.Range("R59").Value = "Fuera de Poligono OVL"
$ws.Range("Q60").Value = "COG-F"
$ws.Range("R60").Value = "Fuera de Poligono OVL"
$ws.Range("Q61").Value = "PUE-J"
$ws.Range("R61").Value = "Fuera de Poligono OVL"
$ws.Range("Q62").Value = "PUE-O"
$ws.Range("R62").Value = "ARATO-25058.PO.1PUE"
$ws.Range("Q63").Value = "COG-O"
$ws.Range("R63").Value = "Fuera de Poligono OVL"
$ws.Range("Q64").Value = "PPT-H"
$ws.Range("R64").Value = "Fuera de Poligono OVL"
$ws.Range("Q65").Value = "PPT-K"
$ws.Range("R65").Value = "Fuera de Poligono OVL"
$ws.Range("Q66").Value = "PPT-O"
$ws.Range("R66").Value = "Fuera de Poligono OVL"
$ws.Range("Q67").Value = "VCR-O"
$ws.Range("R67").Value = "Fuera de Poligono OVL"
$ws.Range("Q68").Value = "COG-C"
$ws.Range("R68").Value = "Fuera de Poligono OVL"
$ws.Range("Q69").Value = "NRA-H"
$ws.Range("R69").Value = "Fuera de Poligono OVL"
$ws.Range("Q70").Value = "PPT-O"
$ws.Range("R70").Value = "Fuera de Poligono OVL"
$ws.Range("Q71").Value = "COG-O"
$ws.Range("R71").Value = "Fuera de Poligono OVL"
$ws.Range("Q72").Value = "DEV-M"
$ws.Range("R72").Value = "ARATO-25058.PO.2DEV"
$ws.Range("Q73").Value = "PAV-?"
$ws.Range("R73").Value = "Fuera de Poligono OVL"
$ws.Range("Q74").Value = "VCR-F"
$ws.Range("R74").Value = "Fuera de Poligono OVL"
$ws.Range("Q75").Value = "CON-M"
$ws.Range("R75").Value = "Fuera de Poligono OVL"
$ws.Range("Q76").Value = "PAV-C"
$ws.Range("R76").Value = "Fuera de Poligono OVL"
$ws.Range("Q77").Value = "PCH-F"
$ws.Range("R77").Value = "Fuera de Poligono OVL"
$ws.Range("Q78").Value = "CLI-D"
$ws.Range("R78").Value = "Fuera de Poligono OVL"
$ws.Range("Q79").Value = "CEN-C"
$ws.Range("R79").Value = "Fuera de Poligono OVL"
$ws.Range("Q80").Value = "CEN-C"
$ws.Range("R80").Value = "Fuera de Poligono OVL"

$ws.Range("A1").Select()
